# Auto-generated Excel COM-interop edit script
# Updates the cryptos list worksheet to match the target commit snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.686.76"
$ws.Range("E2").Value = "  +2.45%  "
$ws.Range("D3").Value = "3.643.42"
$ws.Range("E3").Value = "  +7.02%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "589.64"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").Value = "181.31"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "3.635.82"
$ws.Range("E7").Value = "  +6.98%  "
$ws.Range("D8").Value = "0.618"
$ws.Range("E8").Value = "  +3.21%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "0.202"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("D11").Value = "0.609"
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("D12").Value = "49.84"
$ws.Range("E12").Value = "  +2.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000286"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").Value = "682.18"
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").Value = "4.230.72"
$ws.Range("E15").Value = "  +6.94%  "
$ws.Range("D16").Value = "9.03"
$ws.Range("E16").Value = "  +3.83%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.644.14"
$ws.Range("E17").Value = "  +6.95%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "71.846.64"
$ws.Range("E18").Value = "  +2.65%  "
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").Value = "18.26"
$ws.Range("E20").Value = "  +2.77%  "
$ws.Range("D21").Value = "11.64"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("D22").Value = "0.942"
$ws.Range("E22").Value = "  +2.21%  "
$ws.Range("E23").Value = "  +13.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.90"
$ws.Range("E24").Value = "  +3.18%  "
$ws.Range("D25").Value = "103.21"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("E26").Value = "  +2.20%  "
$ws.Range("D27").Value = "2.83"
$ws.Range("E27").Value = "  +4.05%  "
$ws.Range("D28").Value = "10.12"
$ws.Range("E28").Value = "  +4.59%  "
$ws.Range("D29").Value = "35.33"
$ws.Range("E29").Value = "  +4.24%  "
$ws.Range("D30").Value = "9.21"
$ws.Range("E30").Value = "  +4.05%  "
$ws.Range("D31").Value = "7.35"
$ws.Range("E31").Value = "  +4.74%  "
$ws.Range("D32").Value = "4.19"
$ws.Range("E32").Value = "  +12.87%  "
$ws.Range("D33").Value = "579.28"
$ws.Range("E33").Value = "  +4.20%  "
$ws.Range("D34").Value = "11.34"
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("E35").Value = "  +1.53%  "
$ws.Range("D36").Value = "59.46"
$ws.Range("E36").Value = "  +1.60%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").Value = "3.725.28"
$ws.Range("E38").Value = "  +1.45%  "
$ws.Range("D39").Value = "0.144"
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("D40").Value = "35.68"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "0.0₃0764"
$ws.Range("E41").Value = "  +3.60%  "
$ws.Range("D42").Value = "3.44"
$ws.Range("E42").Value = "  +3.48%  "
$ws.Range("D43").Value = "0.0468"
$ws.Range("E43").Value = "  +9.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.80"
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("D45").Value = "0.346"
$ws.Range("E45").Value = "  +1.89%  "
$ws.Range("D46").Value = "3.42"
$ws.Range("E46").Value = "  +1.70%  "
$ws.Range("E47").Value = "  +5.05%  "
$ws.Range("E48").Value = "  +3.02%  "
$ws.Range("E49").Value = "  +3.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("D51").Value = "132.18"
$ws.Range("E51").Value = "  +1.18%  "
